$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Applications" block used to list App3's five dependencies (rows 37-41)
# right after App2's five dependencies (rows 32-36). The right-hand
# "container" column should instead list things in alphabetical order, so
# App3's row of dependents is replaced by a single "App 2 -> App 3" row
# (continuing the App2 block) and the separate App3 block (rows 38-41) is
# removed entirely.
# ---------------------------------------------------------------------------

# Remove the old App3 dependency rows 38-41 (App3/Procurements, App3/People,
# App3/Data, App3/Technology). This shifts nothing above row 38, and leaves
# what was row 37 (App3/Facilities) as the row to rewrite next.
$ws.Rows("38:41").Delete() | Out-Null

# Rewrite row 37 so the "Applications" block (App2) now also depends on App3
# (alphabetical: App2 -> App3), replacing the old "App3 depends on Facilities"
# row.
$ws.Range("B37").Value2 = "App 2"
$ws.Range("C37").Value2 = "Parent Description…"
$ws.Range("D37").Value2 = "Depends On"
$ws.Range("E37").Value2 = "Applications"
$ws.Range("F37").Value2 = "App 3"
$ws.Range("G37").Value2 = "Dependency Description…"

# A new (otherwise empty) row shows up at the bottom of the used range,
# carrying the header/bold-centered style into column H.
$ws.Range("H42").Font.Bold = $true
$ws.Range("H42").HorizontalAlignment = -4108
$ws.Range("H42").VerticalAlignment = -4108

# Widen the data columns.
$ws.Columns("A").ColumnWidth = 13.233072916666666
$ws.Columns("B").ColumnWidth = 12.764322916666666
$ws.Columns("C").ColumnWidth = 20.963541666666668
$ws.Columns("D").ColumnWidth = 14.166666666666666
$ws.Columns("E").ColumnWidth = 20.764322916666668
$ws.Columns("F").ColumnWidth = 20.631510416666668
$ws.Columns("G").ColumnWidth = 25.166666666666668

# Restore the selection to F5.
$ws.Range("F5").Select() | Out-Null
